$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (columns B,C,D,E,G). Column F is unchanged by this edit.
$data = @{
    2  = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    3  = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    4  = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.997429241610044)
    5  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 14.36450238910742)
    6  = @(0.127881588408715, 0.3127903958511391, 0.8054896365839992, 8.660232485948974, 9.906394106792828)
    7  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    8  = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 3.645393585217082)
    9  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    10 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    11 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    12 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 8.660232485948974, 13.71653804550039)
    13 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    14 = @(0.6753301551942219, 10.29869402782916, 26.21740644021617, 616238.5361209477, 616275.727551571)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G - sum
}
